$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'57.614.35"
$ws.Range("E2").Value = "  -4.08%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'2.926.00"
$ws.Range("E3").Value = "  -2.30%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'549.00"
$ws.Range("E5").Value = "  -4.19%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'130.07"
$ws.Range("E6").Value = "  +3.89%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +1.67%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "'2.919.98"
$ws.Range("E9").Value = "  -2.32%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -3.35%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'4.75"
$ws.Range("E11").Value = "  -6.04%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +1.25%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -0.13%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'32.79"
$ws.Range("E14").Value = "  +0.81%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +0.02%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "'3.404.79"
$ws.Range("E16").Value = "  -2.37%  "

# Row 17 - Polkadot
$ws.Range("D17").Value = "'6.84"
$ws.Range("E17").Value = "  +5.90%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "'2.919.52"
$ws.Range("E18").Value = "  -2.29%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "'57.568.29"
$ws.Range("E19").Value = "  -4.11%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'416.59"
$ws.Range("E20").Value = "  -2.59%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'13.16"
$ws.Range("E21").Value = "  +0.47%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  +1.83%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "'6.96"
$ws.Range("E23").Value = "  -1.26%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("E24").Value = "  +0.64%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "'79.89"
$ws.Range("E25").Value = "  +0.86%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.07%  "

# Row 27 - FirstDigitalUSD
$ws.Range("E27").Value = "  -0.09%  "

# Row 28 - PancakeSwap
$ws.Range("D28").Value = "'2.48"
$ws.Range("E28").Value = "  -2.18%  "

# Row 29 - RenderToken
$ws.Range("D29").Value = "'7.44"
$ws.Range("E29").Value = "  +2.46%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  +1.53%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value = "'25.17"
$ws.Range("E31").Value = "  -0.23%  "

# Row 32 - NEARProtocol
$ws.Range("E32").Value = "  -3.22%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.0972"
$ws.Range("E33").Value = "  +3.34%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  +0.26%  "

# Row 35 - Mantle
$ws.Range("D35").Value = "'0.942"

# Row 36 - Stacks
$ws.Range("E36").Value = "  +0.06%  "

# Row 37 - OKB
$ws.Range("D37").Value = "'47.86"
$ws.Range("E37").Value = "  -4.49%  "

# Row 38 - Cosmos
$ws.Range("D38").Value = "'8.68"
$ws.Range("E38").Value = "  +3.43%  "

# Row 39 - PEPE
$sub3 = [char]0x2083
$ws.Range("D39").Value = "0.0{0}0677" -f $sub3
$ws.Range("E39").Value = "  +2.12%  "

# Row 40 - dogwifhat
$ws.Range("D40").Value = "'2.56"
$ws.Range("E40").Value = "  +3.77%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  -0.79%  "

# Rows 42 & 43 swap: Bittensor/VeChain order flips, with updated values
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0345"
$ws.Range("E42").Value = "  -2.86%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'373.97"
$ws.Range("E43").Value = "  -0.87%  "

# Row 44 - Maker
$ws.Range("D44").Value = "'2.666.48"
$ws.Range("E44").Value = "  -0.08%  "

# Row 45 - USDe
$ws.Range("E45").Value = "  +0.06%  "

# Row 46 - TheGraph
$ws.Range("E46").Value = "  +1.40%  "

# Row 47 - Monero
$ws.Range("D47").Value = "'122.22"
$ws.Range("E47").Value = "  +1.99%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  +1.60%  "

# Row 49 - Fetch.AI
$ws.Range("E49").Value = "  -1.72%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "'23.12"
$ws.Range("E50").Value = "  -1.67%  "

# Row 51 - ThetaToken
$ws.Range("D51").Value = "'1.99"
$ws.Range("E51").Value = "  -0.64%  "
